$wb = $excel.ActiveWorkbook

# --- Insert a new sheet "2022-Q4" right after "总计" (before current "2022-Q3") ---
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")
$newWs = $wb.Worksheets.Add($sheetQ3)
$newWs.Name = "2022-Q4"

# Headers (row 1)
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Data rows (row 2 - 16)
$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "150103"
$newWs.Range("C2").Value = "银河银泰混合"
$newWs.Range("D2").Value = "12.33"
$newWs.Range("E2").Value = "79.30"
$newWs.Range("F2").Value = "4.16"
$newWs.Range("G2").Value = "0.5129"
$newWs.Range("H2").Value = 7

$newWs.Range("A3").Value = 1
$newWs.Range("B3").Value = "519670"
$newWs.Range("C3").Value = "银河行业混合A"
$newWs.Range("D3").Value = "9.81"
$newWs.Range("E3").Value = "92.86"
$newWs.Range("F3").Value = "4.08"
$newWs.Range("G3").Value = "0.4002"
$newWs.Range("H3").Value = 8

$newWs.Range("A4").Value = 2
$newWs.Range("B4").Value = "519679"
$newWs.Range("C4").Value = "银河主题混合"
$newWs.Range("D4").Value = "6.15"
$newWs.Range("E4").Value = "88.50"
$newWs.Range("F4").Value = "4.11"
$newWs.Range("G4").Value = "0.2528"
$newWs.Range("H4").Value = 7

$newWs.Range("A5").Value = 3
$newWs.Range("B5").Value = "014737"
$newWs.Range("C5").Value = "创金合信专精特新股票C"
$newWs.Range("D5").Value = "1.83"
$newWs.Range("E5").Value = "91.82"
$newWs.Range("F5").Value = "6.55"
$newWs.Range("G5").Value = "0.1199"
$newWs.Range("H5").Value = 9

$newWs.Range("A6").Value = 4
$newWs.Range("B6").Value = "014189"
$newWs.Range("C6").Value = "南方专精特新混合A"
$newWs.Range("D6").Value = "2.34"
$newWs.Range("E6").Value = "85.02"
$newWs.Range("F6").Value = "4.55"
$newWs.Range("G6").Value = "0.1065"
$newWs.Range("H6").Value = 6

$newWs.Range("A7").Value = 5
$newWs.Range("B7").Value = "007203"
$newWs.Range("C7").Value = "银河新动能混合"
$newWs.Range("D7").Value = "2.69"
$newWs.Range("E7").Value = "93.01"
$newWs.Range("F7").Value = "3.95"
$newWs.Range("G7").Value = "0.1063"
$newWs.Range("H7").Value = 8

$newWs.Range("A8").Value = 6
$newWs.Range("B8").Value = "011629"
$newWs.Range("C8").Value = "银河核心优势混合A"
$newWs.Range("D8").Value = "2.29"
$newWs.Range("E8").Value = "69.64"
$newWs.Range("F8").Value = "3.99"
$newWs.Range("G8").Value = "0.0914"
$newWs.Range("H8").Value = 8

$newWs.Range("A9").Value = 7
$newWs.Range("B9").Value = "013665"
$newWs.Range("C9").Value = "银河成长优选一年持有混合A"
$newWs.Range("D9").Value = "2.11"
$newWs.Range("E9").Value = "91.96"
$newWs.Range("F9").Value = "3.97"
$newWs.Range("G9").Value = "0.0838"
$newWs.Range("H9").Value = 7

$newWs.Range("A10").Value = 8
$newWs.Range("B10").Value = "519642"
$newWs.Range("C10").Value = "银河大国智造主题灵活配置混合"
$newWs.Range("D10").Value = "1.65"
$newWs.Range("E10").Value = "91.57"
$newWs.Range("F10").Value = "3.96"
$newWs.Range("G10").Value = "0.0653"
$newWs.Range("H10").Value = 7

$newWs.Range("A11").Value = 9
$newWs.Range("B11").Value = "014736"
$newWs.Range("C11").Value = "创金合信专精特新股票A"
$newWs.Range("D11").Value = "0.81"
$newWs.Range("E11").Value = "91.82"
$newWs.Range("F11").Value = "6.55"
$newWs.Range("G11").Value = "0.0531"
$newWs.Range("H11").Value = 9

$newWs.Range("A12").Value = 10
$newWs.Range("B12").Value = "015056"
$newWs.Range("C12").Value = "百嘉百盛混合"
$newWs.Range("D12").Value = "1.17"
$newWs.Range("E12").Value = "82.33"
$newWs.Range("F12").Value = "3.75"
$newWs.Range("G12").Value = "0.0439"
$newWs.Range("H12").Value = 9

$newWs.Range("A13").Value = 11
$newWs.Range("B13").Value = "014190"
$newWs.Range("C13").Value = "南方专精特新混合C"
$newWs.Range("D13").Value = "0.65"
$newWs.Range("E13").Value = "85.02"
$newWs.Range("F13").Value = "4.55"
$newWs.Range("G13").Value = "0.0296"
$newWs.Range("H13").Value = 6

$newWs.Range("A14").Value = 12
$newWs.Range("B14").Value = "013666"
$newWs.Range("C14").Value = "银河成长优选一年持有混合C"
$newWs.Range("D14").Value = "0.60"
$newWs.Range("E14").Value = "91.96"
$newWs.Range("F14").Value = "3.97"
$newWs.Range("G14").Value = "0.0238"
$newWs.Range("H14").Value = 7

$newWs.Range("A15").Value = 13
$newWs.Range("B15").Value = "015670"
$newWs.Range("C15").Value = "银河行业混合C"
$newWs.Range("D15").Value = "0.01"
$newWs.Range("E15").Value = "92.86"
$newWs.Range("F15").Value = "4.08"
$newWs.Range("G15").Value = "0.0004"
$newWs.Range("H15").Value = 8

$newWs.Range("A16").Value = 14
$newWs.Range("B16").Value = "016981"
$newWs.Range("C16").Value = "银河核心优势混合C"
$newWs.Range("D16").Value = "0.00"
$newWs.Range("E16").Value = "69.64"
$newWs.Range("F16").Value = "3.99"
$newWs.Range("G16").Value = 0
$newWs.Range("H16").Value = 8

# --- Update the "总计" summary sheet: insert new 2022-Q4 row, shift others down ---
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 15
$summary.Range("D2").Value = 1.89

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
